$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (Objetivos:) -------------------------------------------------
# B10/C10 long Portuguese objectives text is replaced with the docente text.
$ws.Range("B10").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C10").Value = "198273 - Domingos Savio Giordani"

# --- Row 13 ---------------------------------------------------------------
# Previously row 13 had no A cell (only B13/C13 = "198273 - Domingos Savio Giordani").
# Now it gets a new A13 label "Programa resumido:" (copy column-A formatting from A14
# first, then overwrite the value), and B13/C13 become "Semestral". Row grows to 60pt.
$ws.Range("A14").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14 -----------------------------------------------------------------
# Becomes "Short syllabus:" with the English short-syllabus text.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1 - Basic concepts of Chemistry; 2 - The physical states of matter and their peculiar properties; 3 - Chemical reactions; 4 - Notions of organic chemistry; 5 - Modern materials; 6 - Applied Chemical Technology"
$ws.Range("C14").Value = "1 - Basic concepts of Chemistry; 2 - The physical states of matter and their peculiar properties; 3 - Chemical reactions; 4 - Notions of organic chemistry; 5 - Modern materials; 6 - Applied Chemical Technology"
$ws.Rows(14).RowHeight = 60

# --- Row 15 -------------------------------------------------------------
# Becomes "Programa:" / "01/01/2021". Copy B8:C8 (which already holds the exact text
# "01/01/2021" as a shared string) so that it stays text instead of turning into a date.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8:C8").Copy($ws.Range("B15:C15"))
$ws.Rows(15).RowHeight = 120

# --- Row 16 -----------------------------------------------------------------
# Becomes "Syllabus:" with the long English syllabus text.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1.Basic Concepts of Chemistrya.Atomic Structureb.Periodic tablec.Chemical bonds2.The physical states of matter and their peculiar propertiesa.The gaseous state – pressure, PVT relations, ideal and real gasesb.The liquid state - solutions, intermolecular forces, viscosity, surface tension, vapor pressure, phase changesc.The solid state - classification of solids (molecular, reticular, metallic and ionic)3.Chemical reactionsa.Types of reactions (double-exchange, oxy-reduction)b.Stoichiometry in chemical reactions (limiting reagents, purity and yield) c.Energy and chemical reactionsd.Corrosion Fundamentals4.Notions of organic chemistrya.Hydrocarbons and their main propertiesb.Fuel and combustionc.Polymers5.Applied Chemistry Technologya.Paper and Celluloseb.Sugar and alcoholc.Soaps and detergentsd.Oil and gase.Industrial gasesf.Glass and cement production"
$ws.Range("C16").Value = "1.Basic Concepts of Chemistrya.Atomic Structureb.Periodic tablec.Chemical bonds2.The physical states of matter and their peculiar propertiesa.The gaseous state – pressure, PVT relations, ideal and real gasesb.The liquid state - solutions, intermolecular forces, viscosity, surface tension, vapor pressure, phase changesc.The solid state - classification of solids (molecular, reticular, metallic and ionic)3.Chemical reactionsa.Types of reactions (double-exchange, oxy-reduction)b.Stoichiometry in chemical reactions (limiting reagents, purity and yield) c.Energy and chemical reactionsd.Corrosion Fundamentals4.Notions of organic chemistrya.Hydrocarbons and their main propertiesb.Fuel and combustionc.Polymers5.Applied Chemistry Technologya.Paper and Celluloseb.Sugar and alcoholc.Soaps and detergentsd.Oil and gase.Industrial gasesf.Glass and cement production"
$ws.Rows(16).RowHeight = 120

# --- Row 17 -------------------------------------------------------------
# Becomes just "Avaliação:" with no B/C cells at all (fully cleared, not merely blanked),
# and the row height reverts back to the sheet default.
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows(17).AutoFit()

# --- Row 18 -----------------------------------------------------------------
# Becomes "Método:" / "198273 - Domingos Savio Giordani". B18/C18 did not exist before,
# so copy B10:C10 (already holding the identical text & matching column format) wholesale.
$ws.Range("A18").Value = "Método:"
$ws.Range("B10:C10").Copy($ws.Range("B18:C18"))
$ws.Rows(18).RowHeight = 60

# --- Row 19 -------------------------------------------------------------
# Label changes from "Método:" to "Critério:"; the seminar text in B19/C19 stays the same.
$ws.Range("A19").Value = "Critério:"

# --- Row 20 -------------------------------------------------------------
# Label changes from "Critério:" to "Norma de recuperação:"; B20/C20 text stays the same.
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21 -------------------------------------------------------------
# Label changes from "Norma de recuperação:" to "Bibliografia:"; B21/C21 text stays the
# same, but the row grows from 60pt to 120pt.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# --- Row 22 -------------------------------------------------------------
# Old "Bibliografia:" row (with the long bibliography text) is removed entirely.
$ws.Rows("22:22").Delete()
